$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 29.223446
$ws.Range("H2").Value = 87.670338
$ws.Range("I2").Value = 0.0169041244192178
$ws.Range("J2").Value = 0.0169041244192178
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1591403333333333
$ws.Range("N2").Value = 0.477421
$ws.Range("O2").Value = 0.01260326986877891
$ws.Range("P2").Value = 0.01260326986877891
$ws.Range("Q2").Value = 4.650628937588666
$ws.Range("R2").Value = 41.855660438298
$ws.Range("S2").Value = 0.0002130472419508175
$ws.Range("T2").Value = 0.0002130472419508175
$ws.Range("G3").Value = 29.223446
$ws.Range("H3").Value = 87.670338
$ws.Range("I3").Value = 0.0169041244192178
$ws.Range("J3").Value = 0.0169041244192178
$ws.Range("O3").Value = 0.3005157372251983
$ws.Range("P3").Value = 0.3005157372251983
$ws.Range("Q3").Value = 110.8908400987613
$ws.Range("R3").Value = 998.017560888852
$ws.Range("S3").Value = 0.005079955411987714
$ws.Range("T3").Value = 0.005079955411987715
$ws.Range("G4").Value = 29.223446
$ws.Range("H4").Value = 87.670338
$ws.Range("I4").Value = 0.0169041244192178
$ws.Range("J4").Value = 0.0169041244192178
$ws.Range("M4").Value = 8.673183333333334
$ws.Range("N4").Value = 26.01955
$ws.Range("O4").Value = 0.6868809929060228
$ws.Range("P4").Value = 0.6868809929060229
$ws.Range("Q4").Value = 253.4603047897667
$ws.Range("R4").Value = 2281.1427431079
$ws.Range("S4").Value = 0.01161112176527927
$ws.Range("T4").Value = 0.01161112176527927
$ws.Range("I5").Value = 0.9471112884046843
$ws.Range("J5").Value = 0.9471112884046842
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1591403333333333
$ws.Range("N5").Value = 0.477421
$ws.Range("O5").Value = 0.01260326986877891
$ws.Range("P5").Value = 0.01260326986877891
$ws.Range("Q5").Value = 260.5673654391811
$ws.Range("R5").Value = 2345.10628895263
$ws.Range("S5").Value = 0.01193669916353113
$ws.Range("T5").Value = 0.01193669916353113
$ws.Range("I6").Value = 0.9471112884046843
$ws.Range("J6").Value = 0.9471112884046842
$ws.Range("O6").Value = 0.3005157372251983
$ws.Range("P6").Value = 0.3005157372251983
$ws.Range("R6").Value = 55917.34150213261
$ws.Range("S6").Value = 0.284621847069241
$ws.Range("T6").Value = 0.2846218470692411
$ws.Range("I7").Value = 0.9471112884046843
$ws.Range("J7").Value = 0.9471112884046842
$ws.Range("M7").Value = 8.673183333333334
$ws.Range("N7").Value = 26.01955
$ws.Range("O7").Value = 0.6868809929060228
$ws.Range("P7").Value = 0.6868809929060229
$ws.Range("Q7").Value = 14200.97899634294
$ws.Range("R7").Value = 127808.8109670865
$ws.Range("S7").Value = 0.6505527421719121
$ws.Range("T7").Value = 0.6505527421719121
$ws.Range("G8").Value = 37.39212666666667
$ws.Range("H8").Value = 112.17638
$ws.Range("I8").Value = 0.02162924801792661
$ws.Range("J8").Value = 0.0216292480179266
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1591403333333333
$ws.Range("N8").Value = 0.477421
$ws.Range("O8").Value = 0.01260326986877891
$ws.Range("P8").Value = 0.01260326986877891
$ws.Range("Q8").Value = 5.950595501775556
$ws.Range("R8").Value = 53.55535951598
$ws.Range("S8").Value = 0.0002725992498286804
$ws.Range("T8").Value = 0.0002725992498286804
$ws.Range("G9").Value = 37.39212666666667
$ws.Range("H9").Value = 112.17638
$ws.Range("I9").Value = 0.02162924801792661
$ws.Range("J9").Value = 0.0216292480179266
$ws.Range("O9").Value = 0.3005157372251983
$ws.Range("P9").Value = 0.3005157372251983
$ws.Range("Q9").Value = 141.8875905033911
$ws.Range("R9").Value = 1276.98831453052
$ws.Range("S9").Value = 0.006499929413733872
$ws.Range("T9").Value = 0.006499929413733873
$ws.Range("G10").Value = 37.39212666666667
$ws.Range("H10").Value = 112.17638
$ws.Range("I10").Value = 0.02162924801792661
$ws.Range("J10").Value = 0.0216292480179266
$ws.Range("M10").Value = 8.673183333333334
$ws.Range("N10").Value = 26.01955
$ws.Range("O10").Value = 0.6868809929060228
$ws.Range("P10").Value = 0.6868809929060229
$ws.Range("Q10").Value = 324.3087698032223
$ws.Range("R10").Value = 2918.778928229
$ws.Range("S10").Value = 0.01485671935436405
$ws.Range("T10").Value = 0.01485671935436405
$ws.Range("G11").Value = 24.817167
$ws.Range("H11").Value = 74.45150100000001
$ws.Range("I11").Value = 0.01435533915817136
$ws.Range("J11").Value = 0.01435533915817136
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.1591403333333333
$ws.Range("N11").Value = 0.477421
$ws.Range("O11").Value = 0.01260326986877891
$ws.Range("P11").Value = 0.01260326986877891
$ws.Range("Q11").Value = 3.949412228769
$ws.Range("R11").Value = 35.54471005892101
$ws.Range("S11").Value = 0.0001809242134682831
$ws.Range("T11").Value = 0.0001809242134682831
$ws.Range("G12").Value = 24.817167
$ws.Range("H12").Value = 74.45150100000001
$ws.Range("I12").Value = 0.01435533915817136
$ws.Range("J12").Value = 0.01435533915817136
$ws.Range("O12").Value = 0.3005157372251983
$ws.Range("P12").Value = 0.3005157372251983
$ws.Range("Q12").Value = 94.170841368306
$ws.Range("R12").Value = 847.537572314754
$ws.Range("S12").Value = 0.004314005330235624
$ws.Range("T12").Value = 0.004314005330235624
$ws.Range("G13").Value = 24.817167
$ws.Range("H13").Value = 74.45150100000001
$ws.Range("I13").Value = 0.01435533915817136
$ws.Range("J13").Value = 0.01435533915817136
$ws.Range("M13").Value = 8.673183333333334
$ws.Range("N13").Value = 26.01955
$ws.Range("O13").Value = 0.6868809929060228
$ws.Range("P13").Value = 0.6868809929060229
$ws.Range("Q13").Value = 215.24383920495
$ws.Range("R13").Value = 1937.19455284455
$ws.Range("S13").Value = 0.009860409614467453
$ws.Range("T13").Value = 0.009860409614467455
